$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Fundraising Success" to "Fundraising"
$ws.Name = "Fundraising"

# Add a new "Created" column value (shared string) to column A for the data rows
$ws.Range("A2").Value = "Created"
$ws.Range("A3").Value = "Created"
$ws.Range("A4").Value = "Created"
$ws.Range("A5").Value = "Created"
$ws.Range("A6").Value = "Created"

# Leave the selection where the author last left it
$ws.Range("H17").Select() | Out-Null
